{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph directly above it) that used to sit right\n// after the \"Apostila ou texto...\" bibliography paragraph, while\n// leaving the final blank paragraph (just before the page-break\n// paragraph) untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph (\"Apostila ou texto...\") and the two\n// text paragraphs that must be removed.\nlet anchorIndex = -1;\nlet jupiterIndex = -1;\nlet copyrightIndex = -1;\n\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (anchorIndex === -1 && t.indexOf(\"Apostila ou texto fornecido\") !== -1) {\n    anchorIndex = i;\n  }\n  if (jupiterIndex === -1 && t.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIndex = i;\n  }\n  if (copyrightIndex === -1 && t.indexOf(\"Powered by Jekyll and Github pages\") !== -1) {\n    copyrightIndex = i;\n  }\n}\n\nif (anchorIndex === -1 || jupiterIndex === -1 || copyrightIndex === -1) {\n  throw new Error(\"Could not locate the expected paragraphs to remove.\");\n}\n\n// The blank paragraph that sits between the anchor paragraph and the\n// \"Ver no Jupiter...\" paragraph also needs to go.\nconst toDelete = [];\nfor (let i = anchorIndex + 1; i <= copyrightIndex; i++) {\n  toDelete.push(items[i]);\n}\n\n// Delete from the bottom up so indices/handles of earlier items remain\n// valid while we issue the delete calls.\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n# (and the blank paragraph directly above it) that used to sit right\n# after the \"Apostila ou texto...\" bibliography paragraph, while\n# leaving the final blank paragraph (just before the page-break\n# paragraph) untouched.\n\n$d = $word.ActiveDocument\n\n$anchorIndex = -1\n$jupiterIndex = -1\n$copyrightIndex = -1\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($anchorIndex -eq -1 -and $t -like \"*Apostila ou texto fornecido*\") {\n        $anchorIndex = $i\n    }\n    if ($jupiterIndex -eq -1 -and $t -like \"*Ver no Jupiter*\") {\n        $jupiterIndex = $i\n    }\n    if ($copyrightIndex -eq -1 -and $t -like \"*Powered by Jekyll and Github pages*\") {\n        $copyrightIndex = $i\n    }\n}\n\nif ($anchorIndex -eq -1 -or $jupiterIndex -eq -1 -or $copyrightIndex -eq -1) {\n    throw \"Could not locate the expected paragraphs to remove.\"\n}\n\n# Build one contiguous range spanning from right after the anchor\n# paragraph's own paragraph mark through the end of the copyright\n# paragraph (inclusive of its paragraph mark), then delete it in a\n# single operation so the three in-between paragraphs (the blank one,\n# \"Ver no Jupiter...\", and \"\u00a9 2020...\") disappear together.\n$startRange = $d.Paragraphs.Item($anchorIndex).Range\n$rangeStart = $startRange.End\n$endRange = $d.Paragraphs.Item($copyrightIndex).Range\n$rangeEnd = $endRange.End\n\n$toDelete = $d.Range($rangeStart, $rangeEnd)\n$toDelete.Delete()\n"}
